$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 202 (shifts rows 202:284 down to 203:285,
# and grows the used range from A1:R284 to A1:R285).
$ws.Rows.Item(202).Insert()

# Populate the newly inserted row 202 with the new data point.
$ws.Cells.Item(202, 1).Value = 4
$ws.Cells.Item(202, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(202, 3).Value = "Los Lagos"
$ws.Cells.Item(202, 4).Value = 44900
$ws.Cells.Item(202, 5).Value = 10
$ws.Cells.Item(202, 6).Value = 100112039
$ws.Cells.Item(202, 7).Value = "Ciboulette"
$ws.Cells.Item(202, 8).Value = "Sin especificar"
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 80
$ws.Cells.Item(202, 11).Value = 6000
$ws.Cells.Item(202, 12).Value = 6000
$ws.Cells.Item(202, 13).Value = 6000
$ws.Cells.Item(202, 14).Value = "`$/docena de atados"
$ws.Cells.Item(202, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(202, 16).Value = 2000
$ws.Cells.Item(202, 17).Value = 3
$ws.Cells.Item(202, 18).Value = "Hortaliza"
